# Daily attendance processing - 2025-11-04 11:44:35
# Reorders the "Recorded By" (column G) comma-separated list in each data
# row by moving the last entry to the front (e.g. "System, X" -> "X, System").
# A small set of specific values are left untouched to mirror the exact
# source data used for this processing run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Exact value -> value replacement map, derived from the processing diff.
$map = @{
    "System, system, backup@backdoor.com" = "backup@backdoor.com, System, system";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "System, backup@backdoor.com"         = "backup@backdoor.com, System";
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
